$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.744.76'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.194.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.00%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.43%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.30'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.112'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.426'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.747.32'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.137'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.43%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000172'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.61%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.08%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.835.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.171.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '366.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.520'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.43%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.168'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.71%  '

$ws.Range("E27").Value = '  +0.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0886'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.89'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.33'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.787.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.68'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0697'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.84%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.67'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.91%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0293'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.27%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.39%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.713'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.29%  '

$ws.Range("B45").Value = 'RenzoRestakedETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.236.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.09%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.104'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.976'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.07%  '

$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.76%  '

$ws.Range("B49").Value = 'SuiNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.799'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.04%  '

$ws.Range("E51").Value = '  +0.09%  '
